# 自动更新Excel文件 - 2025-12-21 23:12:39
# Advances the tracking sheet by one day: "today" moves from 2025-12-21 to
# 2025-12-22. For every data row on the active sheet:
#   - column D = total days (总天), column E = remaining days (剩余),
#     column F = start date (开始时间, yyyymmdd integer)
#   - Remaining = TotalDays - (NewToday - StartDate) in days.
#   - If that would drop below 1 (the stock ran out), the shop is
#     restocked today: StartDate resets to NewToday and Remaining resets
#     to TotalDays.

function Get-DaySerial($yy, $mm, $dd) {
    $adj = [int](((14 - $mm) / 12))
    $y2 = $yy + 4800 - $adj
    $m2 = $mm + 12 * $adj - 3
    $jdn = $dd + [int]((153 * $m2 + 2) / 5) + 365 * $y2 + [int]($y2 / 4) - [int]($y2 / 100) + [int]($y2 / 400) - 32045
    return $jdn
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTodaySerial = Get-DaySerial 2025 12 22
$newTodayValue = 20251222

$lastRow = 99
$usedRange = $ws.UsedRange
if ($usedRange -ne $null) {
    $lastRow = $usedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $fRaw = $fCell.Value2

    if ($fRaw -eq $null -or $fRaw -eq "") { continue }

    $fText = [string]$fRaw
    if ($fText.Length -ne 8) { continue }

    $fYear = [int]$fText.Substring(0, 4)
    $fMonth = [int]$fText.Substring(4, 2)
    $fDay = [int]$fText.Substring(6, 2)

    if ($fMonth -lt 1 -or $fMonth -gt 12 -or $fDay -lt 1 -or $fDay -gt 31) { continue }

    $startSerial = Get-DaySerial $fYear $fMonth $fDay

    $dCell = $ws.Cells.Item($r, 4)
    $totalDays = [int]$dCell.Value2

    $elapsed = $newTodaySerial - $startSerial
    $remaining = $totalDays - $elapsed

    $eCell = $ws.Cells.Item($r, 5)

    if ($remaining -lt 1) {
        # Out of stock as of the new day -> restock today.
        $fCell.Value2 = $newTodayValue
        $eCell.Value2 = $totalDays
    } else {
        $eCell.Value2 = $remaining
    }
}
